# Insert two new data rows (new rows 535 and 536) into the "Zapallo italiano"
# price list, pushing the existing rows 535..597 down to 537..599.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at position 535 (shifts old row 535 -> 537, etc.)
$ws.Rows("535:536").Insert()

# ---- New row 535 ----
$ws.Cells.Item(535, 1).Value = 9
$ws.Cells.Item(535, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(535, 3).Value = "Metropolitana"
$ws.Cells.Item(535, 4).Value = 45124
$ws.Cells.Item(535, 5).Value = 13
$ws.Cells.Item(535, 6).Value = 100112032
$ws.Cells.Item(535, 7).Value = "Zapallo italiano"
$ws.Cells.Item(535, 8).Value = "Bola 8"
$ws.Cells.Item(535, 9).Value = "Primera"
$ws.Cells.Item(535, 10).Value = 52
$ws.Cells.Item(535, 11).Value = 13000
$ws.Cells.Item(535, 12).Value = 14000
$ws.Cells.Item(535, 13).Value = 13500
$ws.Cells.Item(535, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(535, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(535, 16).Value = 270
$ws.Cells.Item(535, 17).Value = 50
$ws.Cells.Item(535, 18).Value = "Hortaliza"

# ---- New row 536 ----
$ws.Cells.Item(536, 1).Value = 9
$ws.Cells.Item(536, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(536, 3).Value = "Metropolitana"
$ws.Cells.Item(536, 4).Value = 45124
$ws.Cells.Item(536, 5).Value = 13
$ws.Cells.Item(536, 6).Value = 100112032
$ws.Cells.Item(536, 7).Value = "Zapallo italiano"
$ws.Cells.Item(536, 8).Value = "Sin especificar"
$ws.Cells.Item(536, 9).Value = "Primera"
$ws.Cells.Item(536, 10).Value = 70
$ws.Cells.Item(536, 11).Value = 14000
$ws.Cells.Item(536, 12).Value = 15000
$ws.Cells.Item(536, 13).Value = 14500
$ws.Cells.Item(536, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(536, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(536, 16).Value = 242
$ws.Cells.Item(536, 17).Value = 60
$ws.Cells.Item(536, 18).Value = "Hortaliza"
